$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" conversion text (A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = $ws1.Range("A1").Value()
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 6.55 = 25871.56 pesos"), "✅ 1000 Bs = 6.54 = 25785.51 pesos"
$text = $text -replace [regex]::Escape("✅ 25871.56 pesos = 6.54 = 980.73 Bs"), "✅ 25785.51 pesos = 6.51 = 970.86 Bs"
$ws1.Range("A1").Value = $text

# --- Update "tasas" rate values (N10, O10, N12, O12) ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 152.799
$ws2.Range("O10").Value = 3940
$ws2.Range("N12").Value = 3960
$ws2.Range("O12").Value = 149.1
